$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting so numeric-looking values
# (e.g. "1.00", "68.149.26") are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '68.149.26'
$ws.Range("E2").Value = '  +2.94%  '

# Row 3
$ws.Range("D3").Value = '3.801.38'
$ws.Range("E3").Value = '  +7.06%  '

# Row 4
$ws.Range("E4").Value = '  +0.40%  '

# Row 5
$ws.Range("D5").Value = '418.54'
$ws.Range("E5").Value = '  -0.27%  '

# Row 6
$ws.Range("D6").Value = '138.97'
$ws.Range("E6").Value = '  +4.88%  '

# Row 7
$ws.Range("D7").Value = '3.777.52'
$ws.Range("E7").Value = '  +6.50%  '

# Row 8
$ws.Range("D8").Value = '0.648'
$ws.Range("E8").Value = '  -1.67%  '

# Row 9
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.28%  '

# Row 10
$ws.Range("D10").Value = '0.771'
$ws.Range("E10").Value = '  -1.64%  '

# Row 11
$ws.Range("D11").Value = '0.186'
$ws.Range("E11").Value = '  +10.62%  '

# Row 12
$ws.Range("D12").Value = '0.0000403'
$ws.Range("E12").Value = '  +44.33%  '

# Row 13
$ws.Range("D13").Value = '43.44'
$ws.Range("E13").Value = '  -0.10%  '

# Row 14
$ws.Range("D14").Value = '10.39'
$ws.Range("E14").Value = '  +2.90%  '

# Row 15
$ws.Range("D15").Value = '4.405.30'
$ws.Range("E15").Value = '  +7.66%  '

# Row 16
$ws.Range("E16").Value = '  -0.59%  '

# Row 17
$ws.Range("D17").Value = '3.788.47'
$ws.Range("E17").Value = '  +6.58%  '

# Row 18
$ws.Range("D18").Value = '20.51'
$ws.Range("E18").Value = '  -0.69%  '

# Row 19
$ws.Range("D19").Value = '13.13'
$ws.Range("E19").Value = '  +2.58%  '

# Row 20
$ws.Range("E20").Value = '  +1.48%  '

# Row 21
$ws.Range("D21").Value = '68.191.06'
$ws.Range("E21").Value = '  +3.26%  '

# Row 22
$ws.Range("D22").Value = '444.17'
$ws.Range("E22").Value = '  -2.04%  '

# Row 23
$ws.Range("D23").Value = '15.46'
$ws.Range("E23").Value = '  +17.20%  '

# Row 24
$ws.Range("D24").Value = '89.79'
$ws.Range("E24").Value = '  -0.53%  '

# Row 25
$ws.Range("E25").Value = '  -4.95%  '

# Row 26
$ws.Range("D26").Value = '38.17'
$ws.Range("E26").Value = '  +11.11%  '

# Row 27
$ws.Range("D27").Value = '3.38'
$ws.Range("E27").Value = '  +0.02%  '

# Row 28
$ws.Range("D28").Value = '9.94'
$ws.Range("E28").Value = '  -0.77%  '

# Row 29
$ws.Range("D29").Value = '5.14'
$ws.Range("E29").Value = '  +6.18%  '

# Row 30
$ws.Range("D30").Value = '0.123'
$ws.Range("E30").Value = '  +4.95%  '

# Row 31
$ws.Range("D31").Value = '12.59'
$ws.Range("E31").Value = '  +0.96%  '

# Row 32
$ws.Range("D32").Value = '2.75'
$ws.Range("E32").Value = '  -1.06%  '

# Row 33
$ws.Range("E33").Value = '  -3.11%  '

# Row 34: 'Kaspa' -> 'InjectiveProtocol'
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = '41.65'
$ws.Range("E34").Value = '  +6.54%  '

# Row 35: 'InjectiveProtocol' -> 'Kaspa'
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '0.163'
$ws.Range("E35").Value = '  +0.29%  '

# Row 36
$ws.Range("D36").Value = '57.85'
$ws.Range("E36").Value = '  -0.20%  '

# Row 37
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.12%  '

# Row 38
$ws.Range("D38").Value = '0.0488'
$ws.Range("E38").Value = '  -3.44%  '

# Row 39: 'PEPE' -> 'ThetaToken'
$ws.Range("B39").Value = 'ThetaToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D39").Value = '2.99'
$ws.Range("E39").Value = '  +28.64%  '

# Row 40: 'ThetaToken' -> 'PEPE'
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0711'
$ws.Range("E40").Value = '  -3.53%  '

# Row 41
$ws.Range("D41").Value = '0.147'
$ws.Range("E41").Value = '  -0.12%  '

# Row 42
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.49%  '

# Row 43: 'ApeXProtocol' -> 'LidoDAOToken'
$ws.Range("B43").Value = 'LidoDAOToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D43").Value = '3.39'
$ws.Range("E43").Value = '  +3.80%  '

# Row 44
$ws.Range("D44").Value = '27.26'
$ws.Range("E44").Value = '  +26.48%  '

# Row 45: 'LidoDAOToken' -> 'Monero'
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").Value = '148.05'
$ws.Range("E45").Value = '  +0.68%  '

# Row 46: 'Monero' -> 'ApeXProtocol'
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '3.20'
$ws.Range("E46").Value = '  +24.25%  '

# Row 47
$ws.Range("D47").Value = '2.09'
$ws.Range("E47").Value = '  +4.41%  '

# Row 48
$ws.Range("D48").Value = '2.89'
$ws.Range("E48").Value = '  -6.25%  '

# Row 49
$ws.Range("D49").Value = '2.61'
$ws.Range("E49").Value = '  -5.97%  '

# Row 50
$ws.Range("D50").Value = '4.30'
$ws.Range("E50").Value = '  -4.30%  '

# Row 51
$ws.Range("D51").Value = '0.301'
$ws.Range("E51").Value = '  -3.56%  '
